$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F55").NumberFormat = "@"
$ws.Range("F55").Value = "1885-06-25"
